$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 764.4865192181566
$ws.Range("C2").Value = 418.4496975564101
$ws.Range("D2").Value = 328.1660104757132
$ws.Range("E2").Value = 286.8234418175787
$ws.Range("B3").Value = 871.0289635510992
$ws.Range("C3").Value = 484.8770979451076
$ws.Range("D3").Value = 381.9381445283123
$ws.Range("E3").Value = 337.5397207041048
$ws.Range("B4").Value = 854.9972987449361
$ws.Range("C4").Value = 474.344269817126
$ws.Range("D4").Value = 372.6505164845579
$ws.Range("E4").Value = 329.2111233072053
$ws.Range("B5").Value = 569.1313722990492
$ws.Range("C5").Value = 327.3027464151228
$ws.Range("D5").Value = 254.1870649260764
$ws.Range("E5").Value = 226.7758417148181
$ws.Range("B6").Value = 424.4418026900111
$ws.Range("C6").Value = 268.1912820599164
$ws.Range("D6").Value = 221.405085883514
$ws.Range("E6").Value = 196.612059074901
$ws.Range("B7").Value = 48.45512858836756
$ws.Range("C7").Value = 29.41818687951267
$ws.Range("D7").Value = 23.80803174869665
$ws.Range("E7").Value = 20.8645862999483
$ws.Range("B8").Value = 1685.980398600371
$ws.Range("C8").Value = 1234.037006739788
$ws.Range("D8").Value = 1108.79858268779
$ws.Range("E8").Value = 1076.224503246965
$ws.Range("B9").Value = 584.3834546921156
$ws.Range("C9").Value = 375.1745696129715
$ws.Range("D9").Value = 313.9960287543359
$ws.Range("E9").Value = 281.4538724956348
$ws.Range("B10").Value = 304.4251458185298
$ws.Range("C10").Value = 174.8934120558239
$ws.Range("D10").Value = 139.4316336290295
$ws.Range("E10").Value = 120.7346305157058
$ws.Range("B11").Value = 55.69333481839088
$ws.Range("C11").Value = 30.38056036881715
$ws.Range("D11").Value = 23.994037210493
$ws.Range("E11").Value = 22.41179303283227
$ws.Range("B12").Value = 117.7167208127806
$ws.Range("C12").Value = 75.71151988413907
$ws.Range("D12").Value = 62.5337325697076
$ws.Range("E12").Value = 53.55384386115228
$ws.Range("B13").Value = 163.5233485432103
$ws.Range("C13").Value = 95.50545873374287
$ws.Range("D13").Value = 77.76696643377903
$ws.Range("E13").Value = 68.2463713327362
